$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy styles for column A (index) and column E (match datetime) from the last existing row (97)
$ws.Range("A97").Copy($ws.Range("A98:A103"))
$ws.Range("E97").Copy($ws.Range("E98:E103"))

# Column D holds text that looks numeric ("2023"); force Text format so COM keeps it as a string
$ws.Range("D98:D103").NumberFormat = "@"

# Row 98
$ws.Range("A98").Value = 97
$ws.Range("B98").Value = 'paraguay'
$ws.Range("C98").Value = 'primera-division'
$ws.Range("D98").Value = '2023'
$ws.Range("E98").Value = 45227.02083333334
$ws.Range("F98").Value = 'Tacuary'
$ws.Range("G98").Value = 3
$ws.Range("H98").Value = 'General Caballero JLM'
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 2.59
$ws.Range("K98").Value = '24/10/2023 01:42'
$ws.Range("L98").Value = 2.83
$ws.Range("M98").Value = '28/10/2023 00:21'
$ws.Range("N98").Value = 3.2
$ws.Range("O98").Value = '24/10/2023 01:42'
$ws.Range("P98").Value = 2.99
$ws.Range("Q98").Value = '28/10/2023 00:27'
$ws.Range("R98").Value = 2.91
$ws.Range("S98").Value = '24/10/2023 01:42'
$ws.Range("T98").Value = 2.9
$ws.Range("U98").Value = '28/10/2023 00:27'
$ws.Range("V98").Value = 'https://www.betexplorer.com/football/paraguay/primera-division/tacuary-general-caballero-jlm/UuYRSXoD/'

# Row 99
$ws.Range("A99").Value = 98
$ws.Range("B99").Value = 'paraguay'
$ws.Range("C99").Value = 'primera-division'
$ws.Range("D99").Value = '2023'
$ws.Range("E99").Value = 45227.95833333334
$ws.Range("F99").Value = 'Nacional Asuncion'
$ws.Range("G99").Value = 2
$ws.Range("H99").Value = 'Ameliano'
$ws.Range("I99").Value = 1
$ws.Range("J99").Value = 2.01
$ws.Range("K99").Value = '24/10/2023 01:42'
$ws.Range("L99").Value = 1.9
$ws.Range("M99").Value = '28/10/2023 22:51'
$ws.Range("N99").Value = 3.51
$ws.Range("O99").Value = '24/10/2023 01:42'
$ws.Range("P99").Value = 3.56
$ws.Range("Q99").Value = '28/10/2023 22:56'
$ws.Range("R99").Value = 3.82
$ws.Range("S99").Value = '24/10/2023 01:42'
$ws.Range("T99").Value = 4.39
$ws.Range("U99").Value = '28/10/2023 22:53'
$ws.Range("V99").Value = 'https://www.betexplorer.com/football/paraguay/primera-division/nacional-asuncion-sportivo-ameliano/rqM9iHOD/'

# Row 100
$ws.Range("A100").Value = 99
$ws.Range("B100").Value = 'paraguay'
$ws.Range("C100").Value = 'primera-division'
$ws.Range("D100").Value = '2023'
$ws.Range("E100").Value = 45228.0625
$ws.Range("F100").Value = 'Guarani'
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 'Guairena'
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 1.77
$ws.Range("K100").Value = '25/10/2023 01:48'
$ws.Range("L100").Value = 1.84
$ws.Range("M100").Value = '29/10/2023 01:25'
$ws.Range("N100").Value = 3.5
$ws.Range("O100").Value = '25/10/2023 01:48'
$ws.Range("P100").Value = 3.5
$ws.Range("Q100").Value = '29/10/2023 01:25'
$ws.Range("R100").Value = 4.66
$ws.Range("S100").Value = '25/10/2023 01:48'
$ws.Range("T100").Value = 4.83
$ws.Range("U100").Value = '29/10/2023 01:25'
$ws.Range("V100").Value = 'https://www.betexplorer.com/football/paraguay/primera-division/guarani-guairena-fc/dnQDjyvK/'

# Row 101
$ws.Range("A101").Value = 100
$ws.Range("B101").Value = 'paraguay'
$ws.Range("C101").Value = 'primera-division'
$ws.Range("D101").Value = '2023'
$ws.Range("E101").Value = 45228.875
$ws.Range("F101").Value = 'Olimpia Asuncion'
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 'Cerro Porteno'
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 2.34
$ws.Range("K101").Value = '25/10/2023 02:22'
$ws.Range("L101").Value = 2.64
$ws.Range("M101").Value = '29/10/2023 20:51'
$ws.Range("N101").Value = 3.52
$ws.Range("O101").Value = '25/10/2023 02:22'
$ws.Range("P101").Value = 3.45
$ws.Range("Q101").Value = '29/10/2023 20:50'
$ws.Range("R101").Value = 2.86
$ws.Range("S101").Value = '25/10/2023 02:22'
$ws.Range("T101").Value = 2.74
$ws.Range("U101").Value = '29/10/2023 20:51'
$ws.Range("V101").Value = 'https://www.betexplorer.com/football/paraguay/primera-division/olimpia-asuncion-cerro-porteno/vBFIkegQ/'

# Row 102
$ws.Range("A102").Value = 101
$ws.Range("B102").Value = 'paraguay'
$ws.Range("C102").Value = 'primera-division'
$ws.Range("D102").Value = '2023'
$ws.Range("E102").Value = 45229.91666666666
$ws.Range("F102").Value = 'Sp. Luqueno'
$ws.Range("G102").Value = 2
$ws.Range("H102").Value = 'Sportivo Trinidense'
$ws.Range("I102").Value = 1
$ws.Range("J102").Value = 2.08
$ws.Range("K102").Value = '26/10/2023 00:42'
$ws.Range("L102").Value = 2.15
$ws.Range("M102").Value = '30/10/2023 21:58'
$ws.Range("N102").Value = 3.63
$ws.Range("O102").Value = '26/10/2023 00:42'
$ws.Range("P102").Value = 3.31
$ws.Range("Q102").Value = '30/10/2023 21:52'
$ws.Range("R102").Value = 3.59
$ws.Range("S102").Value = '26/10/2023 00:42'
$ws.Range("T102").Value = 3.75
$ws.Range("U102").Value = '30/10/2023 21:58'
$ws.Range("V102").Value = 'https://www.betexplorer.com/football/paraguay/primera-division/sp-luqueno-sportivo-trinidense/tK1poF1s/'

# Row 103
$ws.Range("A103").Value = 102
$ws.Range("B103").Value = 'paraguay'
$ws.Range("C103").Value = 'primera-division'
$ws.Range("D103").Value = '2023'
$ws.Range("E103").Value = 45230.02083333334
$ws.Range("F103").Value = 'Libertad Asuncion'
$ws.Range("G103").Value = 4
$ws.Range("H103").Value = 'Resistencia'
$ws.Range("I103").Value = 1
$ws.Range("J103").Value = 1.27
$ws.Range("K103").Value = '26/10/2023 00:42'
$ws.Range("L103").Value = 1.23
$ws.Range("M103").Value = '31/10/2023 00:20'
$ws.Range("N103").Value = 6.07
$ws.Range("O103").Value = '26/10/2023 00:42'
$ws.Range("P103").Value = 6.45
$ws.Range("Q103").Value = '31/10/2023 00:20'
$ws.Range("R103").Value = 10.58
$ws.Range("S103").Value = '26/10/2023 00:42'
$ws.Range("T103").Value = 13.27
$ws.Range("U103").Value = '31/10/2023 00:20'
$ws.Range("V103").Value = 'https://www.betexplorer.com/football/paraguay/primera-division/libertad-asuncion-resistencia/lp1lpZGm/'

# Reset column D styling footprint to default (Normal) now that the text values are locked in
$ws.Range("D98:D103").Style = "Normal"

Write-Output "done"
